# v2p14. Compatible with MF-Swift v2212, updated hardpoints.
# Adds a new "Trailer1Axle_f" droplink (AntiRollBar) sheet, based on the
# existing "Bus_Makulu_r" sheet, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# Duplicate the last existing sheet (Bus_Makulu_r) so the new sheet inherits
# the same layout, column widths, styles, tab color and conditional
# formatting as the rest of the Droplink library sheets.
$sourceSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sourceSheet.Copy($null, $sourceSheet)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Trailer1Axle_f"

# Sheet title / instance name shown in H3.
$ws.Range("H3").Value = "Droplink_Trailer1Axle_f"

# Updated hardpoints for the new Trailer1Axle_f droplink.
$ws.Range("F5").Value = 0.05
$ws.Range("G5").Value = 0.6
$ws.Range("H5").Value = 0.19

$ws.Range("F6").Formula = "=0.3-0.15"
$ws.Range("G6").Value = 0.58
$ws.Range("H6").Value = 0.2

$ws.Range("H7").Value = 50
$ws.Range("H8").Value = 0.5

# Make the new sheet the active tab, with H7 selected, same as the source
# commit.
$ws.Activate()
$ws.Range("H7").Select()
